# "Fixed GPS transmission. Temperature in one byte."
#
# This edit reworks the byte layout table on the "packet layout" sheet:
#  - Byte 7 (Temp) now uses a single byte in "Excess 128 format" (a note is
#    added in column D of that row).
#  - GPS longitude/latitude now each take exactly 4 bytes (bytes 8-11 and
#    12-15) instead of being spread out with gaps.
#  - The remaining bytes (PD, EFM, Cloud mean, Cloud st.dev, Rel hum) shift
#    down to close the gap, and a new "Ext temp" byte is appended at the end.
#  - The table now ends at row 26 (was row 28), so the trailing two rows are
#    removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Remove the two trailing rows (old rows 27 & 28) first so the remaining
# rows can be addressed by their final row numbers.
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(27).Delete()

# Row 9 gains a note in column D explaining the temperature encoding,
# styled like the other bold header-ish notes (style index 3 = bold font,
# matching A1:C1's header style).
# This is done before the table below so the new shared strings are
# registered in the same order the original authoring tool produced them
# ("Excess 128 format" before "Ext temp").
$ws.Range("D9").Value = "Excess 128 format"
$ws.Range("D9").Font.Bold = $true

# Column A = Byte index, column B = Contents, column C = Significance.
$rows = @(
    @(9,  7,  "Temp",         0),
    @(10, 8,  "GPS lng",      3),
    @(11, 9,  "GPS lng",      2),
    @(12, 10, "GPS lng",      1),
    @(13, 11, "GPS lng",      0),
    @(14, 12, "GPS lat",      3),
    @(15, 13, "GPS lat",      2),
    @(16, 14, "GPS lat",      1),
    @(17, 15, "GPS lat",      0),
    @(18, 16, "PD",           1),
    @(19, 17, "PD",           0),
    @(20, 18, "EFM",          1),
    @(21, 19, "EFM",          0),
    @(22, 20, "Cloud mean",   1),
    @(23, 21, "Cloud mean",   0),
    @(24, 22, "Cloud st.dev", 0),
    @(25, 23, "Rel hum",      0),
    @(26, 24, "Ext temp",     0)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
}

# Update the active selection to match the new last row of the table.
$ws.Activate()
$ws.Range("D26").Select()

# Cosmetic workbook window size bump.
$excel.ActiveWindow.Height = 14500
